# Generate Report for Handoff
#
# Effect of this edit (derived from the target diff):
#   - File "7fb02bfb-5711-4bfb-8b48-7534ad7c06c3.md" just had a new handoff
#     package generated, so its status flips back from "Ready for handoff"
#     to "In Translation" and its handoff timestamps move forward.
#   - Because of this, on the "Overview" rollup sheet the 7fb02bfb row now
#     sorts ahead of "ebb60103-...md" (rows 3 and 4 swap places).
#   - The zh-cn / de-de detail sheets get the matching Status + Latest
#     Handoff Datetime updates, and their row 3 / row 4 swap the same way.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

# Row 3 becomes the 7fb02bfb entry (status refreshed to "In Translation",
# new handoff-generation timestamp).
$ov.Range("A3").Value = "7fb02bfb-5711-4bfb-8b48-7534ad7c06c3.md"
$ov.Range("B3").Hyperlinks.Item(1).TextToDisplay = "e2e\7fb02bfb-5711-4bfb-8b48-7534ad7c06c3.md"
$ov.Range("C3").Value = ".md"
$ov.Range("D3").Value = ""
$ov.Range("E3").Value = "In Translation"
$ov.Range("F3").Value = "In Translation"
$ov.Range("G3").Value = "2016-11-02 03:55:35"

# Row 4 becomes the ebb60103 entry (unchanged data, just relocated).
$ov.Range("A4").Value = "ebb60103-3a13-4be0-89c9-6678f5550b3a.md"
$ov.Range("B4").Hyperlinks.Item(1).TextToDisplay = "e2e\ebb60103-3a13-4be0-89c9-6678f5550b3a.md"
$ov.Range("C4").Value = ".md"
$ov.Range("D4").Value = ""
$ov.Range("E4").Value = "In Translation"
$ov.Range("F4").Value = "In Translation"
$ov.Range("G4").Value = "2016-11-02 03:47:29"

# ---------------------------------------------------------------------
# zh-cn detail sheet
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

# Row 3 -> 7fb02bfb
$zh.Range("A3").Hyperlinks.Item(1).TextToDisplay = "7fb02bfb-5711-4bfb-8b48-7534ad7c06c3.md"
$zh.Range("C3").Value = "In Translation"
$zh.Range("G3").Value = "7fb02bfb-5711-4bfb-8b48-7534ad7c06c3.4ccae68f039605e25fadd91893142b8cc6ebf79e.zh-cn.xlf"
$zh.Range("H3").Value = "2016-11-02 03:55:22"

# Row 4 -> ebb60103
$zh.Range("A4").Hyperlinks.Item(1).TextToDisplay = "ebb60103-3a13-4be0-89c9-6678f5550b3a.md"
$zh.Range("C4").Value = "In Translation"
$zh.Range("G4").Value = "ebb60103-3a13-4be0-89c9-6678f5550b3a.39a0144d334d68bc06418ee82e42c5e7ed56fa14.zh-cn.xlf"
$zh.Range("H4").Value = "2016-11-02 03:47:18"

# ---------------------------------------------------------------------
# de-de detail sheet
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

# Row 3 -> 7fb02bfb
$de.Range("A3").Hyperlinks.Item(1).TextToDisplay = "7fb02bfb-5711-4bfb-8b48-7534ad7c06c3.md"
$de.Range("C3").Value = "In Translation"
$de.Range("G3").Value = "7fb02bfb-5711-4bfb-8b48-7534ad7c06c3.4ccae68f039605e25fadd91893142b8cc6ebf79e.de-de.xlf"
$de.Range("H3").Value = "2016-11-02 03:55:35"

# Row 4 -> ebb60103
$de.Range("A4").Hyperlinks.Item(1).TextToDisplay = "ebb60103-3a13-4be0-89c9-6678f5550b3a.md"
$de.Range("C4").Value = "In Translation"
$de.Range("G4").Value = "ebb60103-3a13-4be0-89c9-6678f5550b3a.39a0144d334d68bc06418ee82e42c5e7ed56fa14.de-de.xlf"
$de.Range("H4").Value = "2016-11-02 03:47:29"
